$wb = $excel.ActiveWorkbook

# Update the selection on the existing "Calculation" sheet before branching
# off the copy, so its saved cursor position ends up at B14.
$src = $wb.Worksheets.Item("Calculation")
$null = $src.Select()
$null = $src.Range("B14").Select()

# Duplicate the "Calculation" sheet (this also copies its embedded picture)
# and place the copy right after it, naming it "Calculation_new".
$src.Copy($null, $src)
$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "Calculation_new"
$null = $new.Select()

# Update the component values that differ on the new sheet.
$new.Range("B9").Value = 240
$new.Range("B12").Value = 3900
$new.Range("C11").Value = 9100

# B21 used to be a hard-coded number; it becomes a formula on the new sheet.
$new.Range("B21").Formula = "=B14*B11/(B2-B14)"

# B22's formula denominator changes from 10 to 100.
$new.Range("B22").Formula = "=1/(2*PI()*B17*B21/100)"

# Two new helper formulas next to the filter section.
$new.Range("E19").Formula = "=B14*B11/(B2-B14)"
$new.Range("E21").Formula = "=B15*B11/(B15+B11)"

# E22 holds the literal text "=" (not a formula). Build it with a formula
# that evaluates to the text, then convert the result to a plain value so
# it is stored as literal text rather than as a live formula.
$new.Range("E22").Formula = "=CHAR(61)"
$new.Range("E22").Copy()
$new.Range("E22").PasteSpecial(-4163)

$null = $new.Range("B23").Select()
